# Updated cryptos list (prices + 1h volume change) as scraped by the
# GitHub Actions job. Price cells ("D" column) are stored as plain text
# in the workbook (values like "26.308.63" / "1.001" are not real numbers),
# so any value that Excel could misinterpret as a number/date is entered
# with a leading apostrophe to force text entry and keep the exact digits
# (e.g. trailing zeros such as "0.9770") instead of being normalized.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.324.18'
$ws.Range("E2").Value = '  +2.86%  '
$ws.Range("D3").Value = '1.718.08'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''239.41'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '''0.4719'
$ws.Range("E7").Value = '  -1.74%  '
$ws.Range("D8").Value = '''0.2627'
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = '''0.06195'
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").Value = '1.717.33'
$ws.Range("E10").Value = '  +3.00%  '
$ws.Range("D11").Value = '''0.07062'
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = '''15.37'
$ws.Range("E12").Value = '  +3.18%  '
$ws.Range("D13").Value = '''0.5926'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = '''4.401'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '''76.14'
$ws.Range("E15").Value = '  +1.90%  '
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '''1.001'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '26.319.82'
$ws.Range("E18").Value = '  +2.88%  '
$ws.Range("D19").Value = '''0.000006791'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '''11.55'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Value = '1.938.69'
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("D22").Value = '''4.550'
$ws.Range("E22").Value = '  +1.60%  '
$ws.Range("D23").Value = '''8.736'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = '''5.325'
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("D25").Value = '''135.78'
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("D26").Value = '''15.24'
$ws.Range("E26").Value = '  +0.83%  '

# Row 27 (previously BitcoinCash) and row 28 (previously Toncoin) swapped
# ranking positions; Toncoin now ranks 27th and BitcoinCash 28th.
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '''1.407'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '''108.22'
$ws.Range("E28").Value = '  +3.20%  '
$ws.Range("E29").Value = '  +3.33%  '
$ws.Range("D30").Value = '''4.005'
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").Value = '''3.684'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '''0.07727'
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("D33").Value = '''0.04435'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("D34").Value = '''2.616'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '''0.9770'
$ws.Range("E35").Value = '  +2.55%  '
$ws.Range("D36").Value = '''0.6196'
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").Value = '''0.9260'
$ws.Range("E37").Value = '  +6.57%  '
$ws.Range("D38").Value = '''114.29'
$ws.Range("E38").Value = '  +16.78%  '
$ws.Range("D39").Value = '''2.416'
$ws.Range("E39").Value = '  -7.36%  '
$ws.Range("D40").Value = '''1.002'
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = '''1.901'
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").Value = '''0.01482'
$ws.Range("E42").Value = '  -2.51%  '
$ws.Range("D43").Value = '''5.341'
$ws.Range("E43").Value = '  +13.86%  '
$ws.Range("D44").Value = '''0.3814'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("E45").Value = '  +3.12%  '
$ws.Range("D46").Value = '''6.285'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").Value = '''0.05287'
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = '''30.54'
$ws.Range("E48").Value = '  +3.14%  '
$ws.Range("D49").Value = '''7.699'
$ws.Range("E49").Value = '  +3.62%  '
$ws.Range("D50").Value = '''0.3378'
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("D51").Value = '''1.215'
$ws.Range("E51").Value = '  +1.27%  '
